# Atualização automática de JAGUARI.xlsx
$wb = $excel.ActiveWorkbook

# Rename sheets (case-only changes and accent fix)
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the obsolete "Desarquivamentos Pendentes" sheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
